$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# Fix the product name value (add missing hyphen after "248") on both sheets
$wsInput.Range("B1").Value = "248-MS-EI-DB-SAR-REC-NON-RNI-CTRFD-DL-MD-TR-1-ONTIME"
$wsOutput.Range("B1").Value = "248-MS-EI-DB-SAR-REC-NON-RNI-CTRFD-DL-MD-TR-1-ONTIME"

# Update selection on input sheet
$wsInput.Range("B1").Select()

# Make output sheet the active sheet, with selection at B1
$wsOutput.Activate()
$wsOutput.Range("B1").Select()
